$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$src = $ws.Range("E168")
$dst = $ws.Range("E170")
$src.Copy($dst)
$chars = $dst.Characters()
$chars.Text = "https://www.instagram.com/reel/DFpY_7wMXrN/?igsh=ODY2NXduMzl1OGZx"
$chars2 = $dst.Characters()
$chars2.Font.Underline = 1
$chars2.Font.ColorIndex = 11
$chars2.Font.Name = "Calibri"
$chars2.Font.Size = 11
